# Correction in SA algorithm and 746 logs
# Update the Fitness (column C) values for the run_12 log sheet.
# Rows 2-252 hold the per-generation fitness values; after the SA
# algorithm correction these collapse into three plateaus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = 7836
$ws.Range("C8:C26").Value = 7594
$ws.Range("C27:C252").Value = 7569
